$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '36.306.18'
$ws.Range('E2').Value = '  -0.10%  '
$ws.Range('D3').Value = '1.929.42'
$ws.Range('E3').Value = '  -2.54%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').Value = '''240.50'
$ws.Range('E5').Value = '  -1.86%  '
$ws.Range('E6').Value = '  -3.30%  '
$ws.Range('D7').Value = '''1.00'
$ws.Range('E7').Value = '  -0.09%  '
$ws.Range('D8').Value = '''55.94'
$ws.Range('E8').Value = '  -4.93%  '
$ws.Range('D9').Value = '''0.356'
$ws.Range('E9').Value = '  -4.48%  '
$ws.Range('D10').Value = '''0.0827'
$ws.Range('E10').Value = '  -1.12%  '
$ws.Range('E11').Value = '  -2.12%  '
$ws.Range('D12').Value = '2.208.58'
$ws.Range('E12').Value = '  -2.78%  '
$ws.Range('D13').Value = '''20.85'
$ws.Range('E13').Value = '  -9.60%  '
$ws.Range('D14').Value = '''0.793'
$ws.Range('E14').Value = '  -7.52%  '
$ws.Range('D15').Value = '''13.21'
$ws.Range('E15').Value = '  -5.17%  '
$ws.Range('D16').Value = '''5.08'
$ws.Range('E16').Value = '  -6.45%  '
$ws.Range('D17').Value = '1.931.81'
$ws.Range('E17').Value = '  -3.26%  '
$ws.Range('D18').Value = '36.259.24'
$ws.Range('E18').Value = '  +0.13%  '
$ws.Range('B19').Value = 'ShibaInu'
$ws.Range('C19').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D19').Value = '0.0₃0855'
$ws.Range('E19').Value = '  -2.55%  '
$ws.Range('B20').Value = 'Litecoin'
$ws.Range('C20').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D20').Value = '''68.22'
$ws.Range('E20').Value = '  -2.94%  '
$ws.Range('E21').Value = '  -3.54%  '
$ws.Range('D22').Value = '''4.89'
$ws.Range('E22').Value = '  -7.30%  '
$ws.Range('E23').Value = '  -0.02%  '
$ws.Range('D24').Value = '''2.31'
$ws.Range('E24').Value = '  -8.31%  '
$ws.Range('E25').Value = '  -1.58%  '
$ws.Range('D26').Value = '''9.02'
$ws.Range('E26').Value = '  -8.79%  '
$ws.Range('D27').Value = '''160.05'
$ws.Range('E27').Value = '  -1.69%  '
$ws.Range('E28').Value = '  -2.43%  '
$ws.Range('E29').Value = '  -3.85%  '
$ws.Range('E30').Value = '  -2.96%  '
$ws.Range('D31').Value = '''1.09'
$ws.Range('E31').Value = '  -7.03%  '
$ws.Range('D32').Value = '''4.48'
$ws.Range('E32').Value = '  -7.95%  '
$ws.Range('D33').Value = '''0.0613'
$ws.Range('E33').Value = '  -9.41%  '
$ws.Range('B34').Value = 'InternetComputer(DFINITY)'
$ws.Range('C34').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D34').Value = '''4.10'
$ws.Range('E34').Value = '  -6.52%  '
$ws.Range('B35').Value = 'BinanceUSD'
$ws.Range('C35').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D35').Value = '''1.00'
$ws.Range('E35').Value = '  -0.02%  '
$ws.Range('B36').Value = 'WEMIXToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D36').Value = '''1.78'
$ws.Range('E36').Value = '  -1.65%  '
$ws.Range('B37').Value = 'THORChain'
$ws.Range('C37').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D37').Value = '''5.93'
$ws.Range('E37').Value = '  -4.38%  '
$ws.Range('D38').Value = '''2.12'
$ws.Range('E38').Value = '  -5.23%  '
$ws.Range('D39').Value = '''2.98'
$ws.Range('E39').Value = '  +0.66%  '
$ws.Range('D40').Value = '''0.0956'
$ws.Range('E40').Value = '  +0.11%  '
$ws.Range('E41').Value = '  -1.90%  '
$ws.Range('D42').Value = '''0.0207'
$ws.Range('E42').Value = '  -2.80%  '
$ws.Range('D43').Value = '''1.13'
$ws.Range('E43').Value = '  -8.25%  '
$ws.Range('D44').Value = '''15.48'
$ws.Range('E44').Value = '  -3.99%  '
$ws.Range('D45').Value = '1.318.19'
$ws.Range('E45').Value = '  -3.14%  '
$ws.Range('E46').Value = '  -7.58%  '
$ws.Range('D47').Value = '''84.30'
$ws.Range('E47').Value = '  -8.21%  '
$ws.Range('D48').Value = '''6.94'
$ws.Range('E48').Value = '  -6.71%  '
$ws.Range('E49').Value = '  -0.28%  '
$ws.Range('D50').Value = '2.100.91'
$ws.Range('E50').Value = '  -2.69%  '
$ws.Range('D51').Value = '''42.77'
$ws.Range('E51').Value = '  -4.71%  '
